# Edit script: add a "RemoteImage" column (new column I) to the "Geo" sheet,
# containing links to images re-hosted on reportingnotes.com, and update the
# existing "Localimage" column (H) values to use a "www/images/..." prefix
# instead of "images/...".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Geo")

# --- 1. Insert a new column before the old "Attribution" column (I), ---
# --- shifting it (and its data/styles/hyperlink target cell) to J.    ---
$ws.Columns.Item(9).Insert()

# New column header + size (Excel rounds ColumnWidth to whole pixels, so the
# closest achievable stored width is used).
$ws.Range("I1").Value = "RemoteImage"
$ws.Columns.Item(9).ColumnWidth = 31.42

# --- 2. Update "Localimage" (column H) values: images/... -> www/images/... ---
$ws.Range("H2").Value = "www/images/kallunge_carving.jpg"
$ws.Range("H3").Value = "www/images/moraharpa.jpg"
$ws.Range("H4").Value = "www/images/vefsen.png"
$ws.Range("H5").Value = "www/images/esse_harpa.jpg"
$ws.Range("H6").Value = "www/images/kontrabasharpa.jpg"
$ws.Range("H7").Value = "www/images/sigtuna.jpg"
$ws.Range("H8").Value = "www/images/ESI.jpg"
$ws.Range("H9").Value = "www/images/Viola_a_chiavi_Siena_1408.jpg"
$ws.Range("H10").Value = "www/images/tolfta.jpg"
$ws.Range("H11").Value = "www/images/Schluesselfidel.jpg"
$ws.Range("H12").Value = "www/images/Strohfiddel.jpg"
$ws.Range("H13").Value = "www/images/Silverbasharpa.jpg"
$ws.Range("H14").Value = "www/images/lagga-600.jpg"
$ws.Range("H15").Value = "www/images/eric.jpg"
$ws.Range("H16").Value = "www/images/bohlin.jpeg"

# --- 3. Fill in the new "RemoteImage" column (I) with the re-hosted URLs. ---
$ws.Range("I2").Value = "https://reportingnotes.com/wp-content/uploads/2019/05/kallunge_carving.jpg"
$ws.Range("I3").Value = "https://reportingnotes.com/wp-content/uploads/2019/05/moraharpa.jpg"
$ws.Range("I4").Value = "https://reportingnotes.com/wp-content/uploads/2019/05/Vefsen.png"
$ws.Range("I5").Value = "https://reportingnotes.com/wp-content/uploads/2019/05/esse_harpa.jpg"
$ws.Range("I6").Value = "https://reportingnotes.com/wp-content/uploads/2019/05/kontrabasharpa.jpg"
$ws.Range("I7").Value = "https://reportingnotes.com/wp-content/uploads/2019/05/sigtuna.jpg"
$ws.Range("I8").Value = "https://reportingnotes.com/wp-content/uploads/2019/05/ESI.jpg"
$ws.Range("I9").Value = "https://reportingnotes.com/wp-content/uploads/2019/05/Viola_a_chiavi_Siena_1408.jpg"
$ws.Range("I10").Value = "https://reportingnotes.com/wp-content/uploads/2019/05/tolfta.jpg"
$ws.Range("I11").Value = "https://reportingnotes.com/wp-content/uploads/2019/05/Schluesselfidel.jpg"
$ws.Range("I12").Value = "https://reportingnotes.com/wp-content/uploads/2019/05/Strohfiddel.jpg"
$ws.Range("I13").Value = "https://reportingnotes.com/wp-content/uploads/2019/05/Silverbasharpa-1024x365.jpg"
$ws.Range("I14").Value = "https://reportingnotes.com/wp-content/uploads/2019/05/Lagga-600.jpg"
$ws.Range("I15").Value = "https://reportingnotes.com/wp-content/uploads/2019/05/eric.jpg"
$ws.Range("I16").Value = "https://reportingnotes.com/wp-content/uploads/2019/05/Bohlin-1024x704.jpeg"

# --- 4. Re-point the hyperlinks. The column insert shifted the cell values ---
# --- fine, but this COM layer doesn't auto-update existing Hyperlink      ---
# --- objects' anchor ranges, so rebuild the hyperlink collection fully.   ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("G9"), "https://en.wikipedia.org/wiki/Nyckelharpa", "/media/File:Viola_a_chiavi_Siena_1408.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G10"), "https://upload.wikimedia.org/wikipedia/commons/thumb/d/da/Sweden_tolfta_church_angels_with_nyckelharpa.jpg/800px-Sweden_tolfta_church_angels_with_nyckelharpa.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G7"), "https://nyckelharpansforum.net/sigtunanyckel.htm") | Out-Null
$ws.Hyperlinks.Add($ws.Range("J15"), "http://www.ericsahlstrom.se/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G6"), "https://nyckelharpansforum.net/global/piccar/kbhr.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G3"), "http://musikforskning.se/stmonline/vol_9/ternhag/ex1.jpg") | Out-Null

# --- 5. Update the selected cell to match the author's final cursor position. ---
$ws.Range("I5").Select()
